# Update "想去人数" (interested-count) values in the "展览", "演出"
# and "全部类型" sheets to match the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 92
$ws1.Range("F4").Value = 265
$ws1.Range("F6").Value = 24
$ws1.Range("F7").Value = 260
$ws1.Range("F8").Value = 209
$ws1.Range("F9").Value = 1950
$ws1.Range("F10").Value = 347
$ws1.Range("F11").Value = 4573
$ws1.Range("F12").Value = 73
$ws1.Range("F13").Value = 322

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 50

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 92
$ws4.Range("F5").Value = 50
$ws4.Range("F6").Value = 265
$ws4.Range("F8").Value = 24
$ws4.Range("F9").Value = 260
$ws4.Range("F10").Value = 209
$ws4.Range("F13").Value = 1950
$ws4.Range("F14").Value = 347
$ws4.Range("F15").Value = 4573
$ws4.Range("F16").Value = 73
$ws4.Range("F17").Value = 322
